# Insert a new data row before the current row 383 (shifts rows 383-413
# down to 384-414) and populate it with a new Cebolla price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(383).Insert()

$ws.Cells.Item(383, 1).Value = 4
$ws.Cells.Item(383, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(383, 3).Value = "Los Lagos"
$ws.Cells.Item(383, 4).Value = 44578
$ws.Cells.Item(383, 5).Value = 10
$ws.Cells.Item(383, 6).Value = 100112004
$ws.Cells.Item(383, 7).Value = "Cebolla"
$ws.Cells.Item(383, 8).Value = "Sin especificar"
$ws.Cells.Item(383, 9).Value = "Primera"
$ws.Cells.Item(383, 10).Value = 400
$ws.Cells.Item(383, 11).Value = 7000
$ws.Cells.Item(383, 12).Value = 8000
$ws.Cells.Item(383, 13).Value = 7500
$ws.Cells.Item(383, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(383, 15).Value = "Perú"
$ws.Cells.Item(383, 16).Value = 417
$ws.Cells.Item(383, 17).Value = 18
$ws.Cells.Item(383, 18).Value = "Hortaliza"
